$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-record data between row 3 and row 4
# (Date in D, and Volumen/PrecioMinimo/PrecioMaximo/PrecioPromedio in J:M, and Precio $/Kg in P)

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $addr3 = "$col" + "3"
    $addr4 = "$col" + "4"
    $val3 = $ws.Range($addr3).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr3).Value = $val4
    $ws.Range($addr4).Value = $val3
}
